$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2975.3794
$ws.Range("I70").Value = 3442.1904
$ws.Range("J70").Value = 1750
$ws.Range("K70").Value = 10326.5712
$ws.Range("L70").Value = 5250
$ws.Range("M70").Value = -10056.5712
$ws.Range("N70").Value = -5790
$ws.Range("H73").Value = 2975.3794
$ws.Range("I73").Value = 3442.1904
$ws.Range("J73").Value = 1750
$ws.Range("K73").Value = 10326.5712
$ws.Range("L73").Value = 5250
$ws.Range("M73").Value = -9390.5712
$ws.Range("N73").Value = -7122
$ws.Range("H74").Value = 4806
$ws.Range("I74").Value = 4481.1763
$ws.Range("J74").Value = 5496.25
$ws.Range("K74").Value = 4481.1763
$ws.Range("L74").Value = 5496.25
$ws.Range("M74").Value = -3545.1763
$ws.Range("N74").Value = -7368.25
$ws.Range("H76").Value = 17865078
$ws.Range("I76").Value = 33343580
$ws.Range("J76").Value = 5269.231
$ws.Range("K76").Value = 33343580
$ws.Range("L76").Value = 5269.231
$ws.Range("M76").Value = -33343265
$ws.Range("N76").Value = -5899.231
$ws.Range("H77").Value = 4806
$ws.Range("I77").Value = 4481.1763
$ws.Range("J77").Value = 5496.25
$ws.Range("K77").Value = 22405.8815
$ws.Range("L77").Value = 27481.25
$ws.Range("M77").Value = -17725.8815
$ws.Range("N77").Value = -36841.25
$ws.Range("H79").Value = 17865078
$ws.Range("I79").Value = 33343580
$ws.Range("J79").Value = 5269.231
$ws.Range("K79").Value = 33343580
$ws.Range("L79").Value = 5269.231
$ws.Range("M79").Value = -33342488
$ws.Range("N79").Value = -7453.231
$ws.Range("H97").Value = 806.6667
$ws.Range("I97").Value = 290
$ws.Range("J97").Value = 910
$ws.Range("K97").Value = 870
$ws.Range("L97").Value = 2730
$ws.Range("M97").Value = -374
$ws.Range("N97").Value = -3722
$ws.Range("H115").Value = 3136
$ws.Range("J115").Value = 3136
$ws.Range("L115").Value = 9408
$ws.Range("N115").Value = -12542
$ws.Range("H132").Value = 2753.8809
$ws.Range("I132").Value = 1970.7188
$ws.Range("J132").Value = 5260
$ws.Range("K132").Value = 5912.1564
$ws.Range("L132").Value = 15780
$ws.Range("M132").Value = -3382.1564
$ws.Range("N132").Value = -20840
$ws.Range("H137").Value = 6417.8276
$ws.Range("I137").Value = 8312.35
$ws.Range("J137").Value = 2207.7778
$ws.Range("K137").Value = 24937.05
$ws.Range("L137").Value = 6623.3334
$ws.Range("M137").Value = -22387.05
$ws.Range("N137").Value = -11723.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1817138.5
$ws.Range("I32").Value = 2142484
$ws.Range("J32").Value = 4499.143
$ws.Range("K32").Value = 2142484
$ws.Range("L32").Value = 4499.143
$ws.Range("M32").Value = -2142197
$ws.Range("N32").Value = -5073.143
$ws.Range("H61").Value = 956944.5600000001
$ws.Range("I61").Value = 1004622.2
$ws.Range("J61").Value = 913601.25
$ws.Range("K61").Value = 1004622.2
$ws.Range("L61").Value = 913601.25
$ws.Range("M61").Value = -1004410.2
$ws.Range("N61").Value = -914025.25
$ws.Range("H136").Value = 956944.5600000001
$ws.Range("I136").Value = 1004622.2
$ws.Range("J136").Value = 913601.25
$ws.Range("K136").Value = 3013866.6
$ws.Range("L136").Value = 2740803.75
$ws.Range("M136").Value = -3011316.6
$ws.Range("N136").Value = -2745903.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4281.237
$ws.Range("I134").Value = 4420.276
$ws.Range("J134").Value = 3833.2222
$ws.Range("K134").Value = 13260.828
$ws.Range("L134").Value = 11499.6666
$ws.Range("M134").Value = -10725.828
$ws.Range("N134").Value = -16569.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2682.5
$ws.Range("I31").Value = 1562.6786
$ws.Range("J31").Value = 5818
$ws.Range("K31").Value = 1562.6786
$ws.Range("L31").Value = 5818
$ws.Range("M31").Value = -1267.6786
$ws.Range("N31").Value = -6408
$ws.Range("H34").Value = 2682.5
$ws.Range("I34").Value = 1562.6786
$ws.Range("J34").Value = 5818
$ws.Range("K34").Value = 1562.6786
$ws.Range("L34").Value = 5818
$ws.Range("M34").Value = -1360.6786
$ws.Range("N34").Value = -6222
$ws.Range("H99").Value = 58936.777
$ws.Range("I99").Value = 169043.67
$ws.Range("J99").Value = 3883.3333
$ws.Range("K99").Value = 169043.67
$ws.Range("L99").Value = 3883.3333
$ws.Range("M99").Value = -167545.67
$ws.Range("N99").Value = -6879.3333
$ws.Range("H126").Value = 58936.777
$ws.Range("I126").Value = 169043.67
$ws.Range("J126").Value = 3883.3333
$ws.Range("K126").Value = 507131.01
$ws.Range("L126").Value = 11649.9999
$ws.Range("M126").Value = -504661.01
$ws.Range("N126").Value = -16589.9999
$ws.Range("H132").Value = 1987.3889
$ws.Range("I132").Value = 1016.5926
$ws.Range("J132").Value = 4899.778
$ws.Range("K132").Value = 3049.7778
$ws.Range("L132").Value = 14699.334
$ws.Range("M132").Value = -519.7777999999998
$ws.Range("N132").Value = -19759.334
$ws.Range("H134").Value = 2221.6667
$ws.Range("I134").Value = 1234.5555
$ws.Range("J134").Value = 3702.3333
$ws.Range("K134").Value = 3703.6665
$ws.Range("L134").Value = 11106.9999
$ws.Range("M134").Value = -1168.6665
$ws.Range("N134").Value = -16176.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18852.5
$ws.Range("I70").Value = 4450
$ws.Range("J70").Value = 20910
$ws.Range("K70").Value = 4450
$ws.Range("L70").Value = 20910
$ws.Range("M70").Value = -4180
$ws.Range("N70").Value = -21450
$ws.Range("H73").Value = 18852.5
$ws.Range("I73").Value = 4450
$ws.Range("J73").Value = 20910
$ws.Range("K73").Value = 4450
$ws.Range("L73").Value = 20910
$ws.Range("M73").Value = -3514
$ws.Range("N73").Value = -22782
$ws.Range("H80").Value = 5756.087
$ws.Range("I80").Value = 6974.375
$ws.Range("J80").Value = 2971.4285
$ws.Range("K80").Value = 6974.375
$ws.Range("L80").Value = 2971.4285
$ws.Range("M80").Value = -5976.375
$ws.Range("N80").Value = -4967.4285
$ws.Range("H83").Value = 5756.087
$ws.Range("I83").Value = 6974.375
$ws.Range("J83").Value = 2971.4285
$ws.Range("K83").Value = 34871.875
$ws.Range("L83").Value = 14857.1425
$ws.Range("M83").Value = -29879.875
$ws.Range("N83").Value = -24841.1425
$ws.Range("H132").Value = 4757.9287
$ws.Range("I132").Value = 4548.6313
$ws.Range("J132").Value = 5199.778
$ws.Range("K132").Value = 13645.8939
$ws.Range("L132").Value = 15599.334
$ws.Range("M132").Value = -11115.8939
$ws.Range("N132").Value = -20659.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 560.04
$ws.Range("I22").Value = 457.6
$ws.Range("J22").Value = 628.3333
$ws.Range("K22").Value = 457.6
$ws.Range("L22").Value = 628.3333
$ws.Range("M22").Value = -162.6
$ws.Range("N22").Value = -1218.3333
$ws.Range("H27").Value = 560.04
$ws.Range("I27").Value = 457.6
$ws.Range("J27").Value = 628.3333
$ws.Range("K27").Value = 457.6
$ws.Range("L27").Value = 628.3333
$ws.Range("M27").Value = -350.6
$ws.Range("N27").Value = -842.3333
$ws.Range("H61").Value = 997
$ws.Range("I61").Value = 1025
$ws.Range("J61").Value = 941
$ws.Range("K61").Value = 1025
$ws.Range("L61").Value = 941
$ws.Range("M61").Value = -823
$ws.Range("N61").Value = -1345
$ws.Range("H82").Value = 1732.5834
$ws.Range("I82").Value = 1186
$ws.Range("J82").Value = 2825.75
$ws.Range("K82").Value = 1186
$ws.Range("L82").Value = 2825.75
$ws.Range("M82").Value = -825
$ws.Range("N82").Value = -3547.75
$ws.Range("H85").Value = 1732.5834
$ws.Range("I85").Value = 1186
$ws.Range("J85").Value = 2825.75
$ws.Range("K85").Value = 1186
$ws.Range("L85").Value = 2825.75
$ws.Range("M85").Value = 62
$ws.Range("N85").Value = -5321.75
$ws.Range("H113").Value = 997
$ws.Range("I113").Value = 1025
$ws.Range("J113").Value = 941
$ws.Range("K113").Value = 1025
$ws.Range("L113").Value = 941
$ws.Range("M113").Value = 1145
$ws.Range("N113").Value = -5281
$ws.Range("H132").Value = 11117728
$ws.Range("I132").Value = 15627058
$ws.Range("K132").Value = 46881174
$ws.Range("M132").Value = -46878644

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 930
$ws.Range("I126").Value = 676.6667
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 2030.0001
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = 439.9999
$ws.Range("N126").Value = -9440
